# Apply updated "Inscritos"/"Pagos"/"Inscrições homologadas" counts
# to the Resumo de Inscrições worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value
$updates = @{
    5  = @{ E = 130 }
    10 = @{ E = 456 }
    11 = @{ E = 315 }
    12 = @{ E = 451; F = 246; H = 246 }
    15 = @{ E = 147; F = 62;  H = 62 }
    23 = @{ E = 185 }
    24 = @{ E = 194 }
    25 = @{ E = 244; F = 117; H = 117 }
    26 = @{ E = 145; F = 88;  H = 88 }
    27 = @{ E = 305; F = 147; H = 147 }
    28 = @{ E = 185 }
    29 = @{ E = 158 }
    30 = @{ E = 194 }
    33 = @{ E = 265 }
    37 = @{ E = 144 }
    39 = @{ E = 169 }
    41 = @{ E = 370 }
    42 = @{ E = 338 }
    44 = @{ E = 295; F = 147; H = 147 }
    46 = @{ E = 288; F = 159; H = 159 }
    47 = @{ E = 415 }
    48 = @{ E = 188 }
    49 = @{ E = 269 }
    50 = @{ E = 229; F = 105; H = 105 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}

$wb.Save()
